$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the NmerMatch method label to credit J Greenbaum
$ws.Range("A3").Value = "NmerMatch (J Greenbaum)"

# Update benchmarking numbers for the non-text-shifting methods
$ws.Range("B2").Value = 39.617
$ws.Range("D2").Value = 0.081
$ws.Range("E2").Value = 39.697

$ws.Range("B3").Value = 53.716
$ws.Range("C3").Value = 0.0055389404296875
$ws.Range("D3").Value = 12.319
$ws.Range("E3").Value = 66.041

$ws.Range("B8").Value = 1.265
$ws.Range("D8").Value = 11.29
$ws.Range("E8").Value = 12.555

$ws.Range("B9").Value = 0.25
$ws.Range("D9").Value = 5.011
$ws.Range("E9").Value = 5.262

$ws.Range("B10").Value = 2.645
$ws.Range("D10").Value = 0.497
$ws.Range("E10").Value = 3.142

# Move the active selection to the last updated row's data range
$ws.Range("B10:G10").Select()
